$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.714.98"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.600.47"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'211.25"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'0.511"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").Value = "'0.246"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "1.825.19"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "1.584.68"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "'65.12"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "26.687.02"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").Value = "'210.72"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("E20").Value = "  +2.64%  "
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").Value = "'2.29"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").Value = "'143.84"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").Value = "'0.0514"
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("D34").Value = "1.295.27"
$ws.Range("E34").Value = "  +1.78%  "
$ws.Range("D35").Value = "'2.47"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'1.50"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.605"
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("E38").Value = "  +19.02%  "
$ws.Range("D39").Value = "'0.0169"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "'0.822"
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").Value = "'0.781"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").Value = "'63.23"
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("D45").Value = "1.737.37"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Value = "'90.86"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("D47").Value = "'1.56"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  -1.79%  "
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("E51").Value = "  +0.13%  "
